$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new conversation rows (123-133) to match the exported log snapshot.
# Column D sometimes holds a numeric-looking phone number; Excel would otherwise
# auto-coerce that to a Number cell, so for those we briefly force a text format,
# assign the value, then strip the format again so no stray style is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 123
$ws.Range("A123").Value = '2026-01-10 20:06:32'
$ws.Range("B123").Value = 'Noah Dubitzky'
$ws.Range("C123").Value = 8450689526
Set-TextValue $ws.Range("D123") '13052054965'
$ws.Range("E123").Value = 'Test'

# Row 124
$ws.Range("A124").Value = '2026-01-10 20:36:50'
$ws.Range("B124").Value = 'Noah Dubitzky'
$ws.Range("C124").Value = 8450689526
Set-TextValue $ws.Range("D124") '13052054965'
$ws.Range("E124").Value = 'Test'

# Row 125
$ws.Range("A125").Value = '2026-01-10 20:37:37'
$ws.Range("B125").Value = 'Noah Dubitzky'
$ws.Range("C125").Value = 8450689526
Set-TextValue $ws.Range("D125") '13052054965'
$ws.Range("E125").Value = 'Test'

# Row 126
$ws.Range("A126").Value = '2026-01-10 20:37:44'
$ws.Range("B126").Value = 'Noah Dubitzky'
$ws.Range("C126").Value = 8450689526
Set-TextValue $ws.Range("D126") '13052054965'
$ws.Range("E126").Value = 'Test'

# Row 127
$ws.Range("A127").Value = '2026-01-09 21:30:01'
$ws.Range("B127").Value = 'Noahs life'
$ws.Range("D127").Value = 'Unknown'
$ws.Range("E127").Value = 'Test'
$ws.Range("G127").Value = 'Noahs life'

# Row 128
$ws.Range("A128").Value = '2026-01-10 20:38:50'
$ws.Range("B128").Value = 'Noahs life'
$ws.Range("D128").Value = 'Unknown'
$ws.Range("E128").Value = 'Test'
$ws.Range("G128").Value = 'Noahs life'

# Row 129
$ws.Range("A129").Value = '2026-01-10 20:40:05'
$ws.Range("B129").Value = 'Noah Dubitzky'
$ws.Range("C129").Value = 8450689526
Set-TextValue $ws.Range("D129") '13052054965'
$ws.Range("E129").Value = 'Idk'

# Row 130
$ws.Range("A130").Value = '2026-01-10 20:40:52'
$ws.Range("B130").Value = 'Emerson Walker'
$ws.Range("C130").Value = 8483444103
$ws.Range("D130").Value = 'Unknown'
$ws.Range("E130").Value = 'My name is Emerson and its a please to talk to you.'

# Row 131
$ws.Range("A131").Value = '2026-01-10 20:41:22'
$ws.Range("B131").Value = 'Emerson Walker'
$ws.Range("C131").Value = 8483444103
$ws.Range("D131").Value = 'Unknown'
$ws.Range("E131").Value = 'I’d like to take a moment to explain the purpose of my outreach.'

# Row 132
$ws.Range("A132").Value = '2026-01-10 20:38:50'
$ws.Range("B132").Value = 'Noahs life'
$ws.Range("D132").Value = 'Unknown'
$ws.Range("E132").Value = 'Test'
$ws.Range("G132").Value = 'Noahs life'

# Row 133
$ws.Range("A133").Value = '2026-01-10 20:45:44'
$ws.Range("B133").Value = 'Noahs life'
$ws.Range("D133").Value = 'Unknown'
$ws.Range("E133").Value = 'Test'
$ws.Range("G133").Value = 'Noahs life'

